$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look numeric so Excel keeps them as text
# (matches source data which stores these as literal strings, e.g. "0.510").
# NOTE: applied per-cell (not as a unioned multi-area range) for reliability.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '56.492.98'
$ws.Range('E2').Value = '  +1.63%  '

# Row 3
$ws.Range('D3').Value = '2.480.64'
$ws.Range('E3').Value = '  -0.88%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').Value = '490.13'
$ws.Range('E5').Value = '  +0.81%  '

# Row 6
$ws.Range('D6').Value = '151.22'
$ws.Range('E6').Value = '  +7.27%  '

# Row 7
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.27%  '

# Row 8
$ws.Range('D8').Value = '0.510'
$ws.Range('E8').Value = '  -0.37%  '

# Row 9
$ws.Range('D9').Value = '2.492.32'
$ws.Range('E9').Value = '  -0.37%  '

# Row 10
$ws.Range('D10').Value = '5.73'
$ws.Range('E10').Value = '  +3.58%  '

# Row 11
$ws.Range('D11').Value = '0.0982'
$ws.Range('E11').Value = '  -0.66%  '

# Row 12
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  +0.89%  '

# Row 13
$ws.Range('E13').Value = '  +0.69%  '

# Row 14
$ws.Range('D14').Value = '2.921.02'
$ws.Range('E14').Value = '  -0.77%  '

# Row 15
$ws.Range('D15').Value = '56.708.67'
$ws.Range('E15').Value = '  +1.88%  '

# Row 16
$ws.Range('D16').Value = '21.14'
$ws.Range('E16').Value = '  +1.71%  '

# Row 17
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  -1.70%  '

# Row 18
$ws.Range('D18').Value = '2.491.52'
$ws.Range('E18').Value = '  -0.75%  '

# Row 19
$ws.Range('D19').Value = '4.54'
$ws.Range('E19').Value = '  +3.29%  '

# Row 20
$ws.Range('D20').Value = '10.26'
$ws.Range('E20').Value = '  +2.52%  '

# Row 21
$ws.Range('D21').Value = '320.00'
$ws.Range('E21').Value = '  -0.95%  '

# Row 22
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.27%  '

# Row 23
$ws.Range('D23').Value = '5.86'
$ws.Range('E23').Value = '  +1.86%  '

# Row 24
$ws.Range('D24').Value = '58.32'

# Row 25
$ws.Range('D25').Value = '0.409'
$ws.Range('E25').Value = '  -0.80%  '

# Row 26
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.11%  '

# Row 27
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -4.83%  '

# Row 28
$ws.Range('D28').Value = '2.599.19'
$ws.Range('E28').Value = '  -1.08%  '

# Row 29
$ws.Range('D29').Value = '7.56'
$ws.Range('E29').Value = '  +1.33%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0798'
$ws.Range('E30').Value = '  +0.17%  '

# Row 31
$ws.Range('E31').Value = '  -0.12%  '

# Row 32
$ws.Range('D32').Value = '151.34'
$ws.Range('E32').Value = '  +0.66%  '

# Row 33
$ws.Range('D33').Value = '18.27'
$ws.Range('E33').Value = '  -0.09%  '

# Row 34
$ws.Range('E34').Value = '  +1.06%  '

# Row 35
$ws.Range('D35').Value = '5.25'
$ws.Range('E35').Value = '  +0.34%  '

# Row 36
$ws.Range('D36').Value = '1.16'
$ws.Range('E36').Value = '  +3.53%  '

# Row 37
$ws.Range('D37').Value = '3.74'
$ws.Range('E37').Value = '  +0.88%  '

# Row 38
$ws.Range('D38').Value = '0.868'
$ws.Range('E38').Value = '  -0.64%  '

# Row 39
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '34.21'
$ws.Range('E39').Value = '  -0.25%  '

# Row 40
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.38'
$ws.Range('E40').Value = '  +4.65%  '

# Row 41
$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  +1.92%  '

# Row 42
$ws.Range('D42').Value = '0.0561'
$ws.Range('E42').Value = '  +0.52%  '

# Row 43
$ws.Range('D43').Value = '0.613'
$ws.Range('E43').Value = '  -0.41%  '

# Row 44
$ws.Range('E44').Value = '  -0.27%  '

# Row 45
$ws.Range('D45').Value = '267.07'
$ws.Range('E45').Value = '  +4.63%  '

# Row 46
$ws.Range('D46').Value = '4.86'
$ws.Range('E46').Value = '  +3.32%  '

# Row 47
$ws.Range('D47').Value = '0.0927'
$ws.Range('E47').Value = '  +1.42%  '

# Row 48
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').Value = '10.22'
$ws.Range('E48').Value = '  +1.03%  '

# Row 49
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0229'
$ws.Range('E49').Value = '  +1.29%  '

# Row 50
$ws.Range('D50').Value = '17.76'
$ws.Range('E50').Value = '  +1.18%  '

# Row 51
$ws.Range('D51').Value = '1.878.74'
$ws.Range('E51').Value = '  -6.18%  '
